# Apply the commit: "corrected the results with updated data of changing
# inclusion levels to 20% of Sarson, Mott Grss and Barseem"
#
# 1. Delete the "Translation" sheet (no longer present in target workbook).
# 2. Rename "NDF Added with updated Inclusio" -> "Composition April 6, 2023".
# 3. Update the inclusion-level column (I) on that sheet:
#       Barseem (row 2)          0.6 -> 0.2
#       Maize (row 3)            0.6 -> 0.4
#       Oat (Jai) (row 4)        0.6 -> 0.4
#       Mustard (Sarson) (row 5) 0.4 -> 0.2
#       Maize Silage (row 6)     0.6 -> 0.5
#       Sugarcane (row 7)        0.5 -> 0.4
#       Mott grass (row 9)       0.4 -> 0.2
# 4. Make the renamed sheet the active / selected sheet.

$wb = $excel.ActiveWorkbook

# --- delete the "Translation" sheet -----------------------------------
$wb.Worksheets.Item("Translation").Delete()

# --- rename the NDF sheet ----------------------------------------------
$ndf = $wb.Worksheets.Item("NDF Added with updated Inclusio")
$ndf.Name = "Composition April 6, 2023"

# --- update the inclusion-level values ---------------------------------
$ndf.Range("I2").Value = 0.2
$ndf.Range("I3").Value = 0.4
$ndf.Range("I4").Value = 0.4
$ndf.Range("I5").Value = 0.2
$ndf.Range("I6").Value = 0.5
$ndf.Range("I7").Value = 0.4
$ndf.Range("I9").Value = 0.2

# --- make this sheet the active sheet -----------------------------------
$ndf.Activate()
$ndf.Range("I10").Select()
